$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values (column D) - forced to text to avoid numeric auto-conversion,
# then formats cleared so the cell keeps its original (default) style.
$priceUpdates = @{
    'D2' = '29.979.90'
    'D3' = '1.907.55'
    'D5' = '319.41'
    'D6' = '1.001'
    'D7' = '0.5034'
    'D8' = '0.4039'
    'D9' = '0.08270'
    'D10' = '41.94'
    'D11' = '1.098'
    'D12' = '24.25'
    'D13' = '1.909.39'
    'D14' = '6.382'
    'D15' = '7.208'
    'D16' = '1.004'
    'D17' = '91.82'
    'D18' = '0.00001093'
    'D19' = '0.06513'
    'D20' = '18.04'
    'D21' = '1.001'
    'D22' = '5.933'
    'D23' = '30.011.04'
    'D24' = '11.25'
    'D25' = '2.198'
    'D26' = '22.24'
    'D27' = '2.128.97'
    'D28' = '161.55'
    'D29' = '2.274'
    'D30' = '128.59'
    'D31' = '1.123'
    'D32' = '0.1033'
    'D33' = '5.927'
    'D34' = '3.793'
    'D35' = '5.386'
    'D36' = '0.02433'
    'D37' = '0.06337'
    'D38' = '0.2143'
    'D39' = '0.6509'
    'D40' = '1.193'
    'D41' = '8.638'
    'D42' = '11.33'
    'D43' = '1.207'
    'D44' = '2.209'
    'D45' = '13.30'
    'D46' = '0.6016'
    'D47' = '3.631'
    'D48' = '122.69'
    'D49' = '1.206'
    'D50' = '78.30'
    'D51' = '1.132'
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# Volume/percentage values (column E) - stored as plain text naturally.
$volumeUpdates = @{
    'E2' = '  -0.93%  '
    'E3' = '  +0.06%  '
    'E4' = '  +0.06%  '
    'E5' = '  -2.04%  '
    'E6' = '  +0.09%  '
    'E7' = '  -2.49%  '
    'E8' = '  +0.32%  '
    'E9' = '  -2.53%  '
    'E10' = '  -1.86%  '
    'E11' = '  -1.72%  '
    'E12' = '  +2.81%  '
    'E13' = '  +0.07%  '
    'E14' = '  -1.14%  '
    'E15' = '  -1.90%  '
    'E16' = '  +0.30%  '
    'E18' = '  -1.80%  '
    'E19' = '  -2.38%  '
    'E20' = '  -1.58%  '
    'E21' = '  +0.09%  '
    'E22' = '  -0.93%  '
    'E23' = '  -0.84%  '
    'E24' = '  -0.29%  '
    'E25' = '  -1.22%  '
    'E26' = '  +2.14%  '
    'E27' = '  +0.11%  '
    'E28' = '  +0.13%  '
    'E29' = '  -5.19%  '
    'E30' = '  -0.69%  '
    'E31' = '  +2.24%  '
    'E32' = '  -2.33%  '
    'E33' = '  -1.67%  '
    'E34' = '  +0.62%  '
    'E35' = '  +2.46%  '
    'E36' = '  -2.75%  '
    'E37' = '  -3.73%  '
    'E38' = '  -3.33%  '
    'E39' = '  -0.07%  '
    'E40' = '  -3.42%  '
    'E41' = '  -2.05%  '
    'E42' = '  -4.95%  '
    'E43' = '  -2.39%  '
    'E44' = '  +7.29%  '
    'E45' = '  +0.19%  '
    'E46' = '  -1.86%  '
    'E47' = '  -2.24%  '
    'E48' = '  -1.94%  '
    'E49' = '  -3.07%  '
    'E50' = '  -1.25%  '
    'E51' = '  -2.41%  '
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
